$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("json_path_after_change" / C)
# so the table becomes: A=json_path_before_change, B=json_path_after_change,
# C=filtre_visuels (new), D=instructions (was C)
$ws.Columns("C").Insert()

# Remove the last four rows (the "section" entries), which were rows 12-15
$ws.Rows("12:15").Delete()

# Header for the new column
$ws.Range("C1").Value = "filtre_visuels"

# Fill the new column for every data row (2-11) with the same value
$ws.Range("C2:C11").Value = "slicer, advancedSlicerVisual"

# Match the new column's width to the neighbouring columns
$ws.Columns("C").ColumnWidth = 47.666666666666664

# Update the active selection to match the new used range
$ws.Range("C7").Select()
